$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per repulled data / mean calculation
$updates = @{
    4  = -4
    5  = -3
    6  = -5
    8  = -5
    12 = 2
    15 = -4
    16 = -4
    20 = -3
    23 = -1
    25 = 6
    29 = 2
    33 = -9
    35 = -3
    37 = -2
    39 = -4
    40 = 2
    41 = -1
    45 = 5
    49 = -2
    50 = -2
    52 = 3
    53 = -1
    57 = 0
    60 = 0
    64 = 2
    65 = -3
    69 = 3
    74 = -6
    76 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
